# Updates the cryptos list (prices and 1h volume % changes) on Sheet1.
# Values are assigned with a leading apostrophe so Excel always stores them
# as literal text (matching the original inlineStr cells), never as numbers,
# dates, or other auto-detected types. The Style reset avoids leaving a
# "Text" number-format style attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.068.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.268.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.92%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'586.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.79%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'184.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +4.06%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.89%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.836.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.85%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.45%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +2.67%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.058.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.36%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.267.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'382.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.59%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.44%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +6.70%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.65%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.43%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.37%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'USDe"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'163.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.69%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.837"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'26.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Filecoin"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'4.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.77%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'dogwifhat"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'2.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.71%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'25.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.61%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'41.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.97%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.627.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.90%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'340.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.32%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +2.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'32.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.72%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.21%  "
$ws.Range("E51").Style = "Normal"
